$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.05053916972009631
$ws.Range("D2").Value = 0.001474349248763218
$ws.Range("E2").Value = 0.05083759292226131
$ws.Range("F2").Value = 0.5092771424254181
$ws.Range("G2").Value = 0.3939498039479759
$ws.Range("H2").Value = 0.4304703099388547
$ws.Range("M2").Value = 0.8886877327032181
$ws.Range("N2").Value = 1.779252987943323
$ws.Range("O2").Value = 1.589036467503917
$ws.Range("C3").Value = 0.04480434375342668
$ws.Range("D3").Value = 0.001283251391759777
$ws.Range("E3").Value = 0.05411985105376549
$ws.Range("F3").Value = 0.4792200641138038
$ws.Range("G3").Value = 0.362067737511552
$ws.Range("H3").Value = 0.4199052581246576
$ws.Range("M3").Value = 0.7810507139676162
$ws.Range("N3").Value = 1.628484058654209
$ws.Range("O3").Value = 1.498953055692624
$ws.Range("C4").Value = 0.0412986163176754
$ws.Range("D4").Value = 0.001165518968425872
$ws.Range("E4").Value = 0.05628346368258419
$ws.Range("F4").Value = 0.4611257301856426
$ws.Range("G4").Value = 0.3427613838968853
$ws.Range("H4").Value = 0.4137109623053874
$ws.Range("M4").Value = 0.7148008982813536
$ws.Range("N4").Value = 1.536058465213614
$ws.Range("O4").Value = 1.444810740166986
$ws.Range("C5").Value = 0.03987384526845972
$ws.Range("D5").Value = 0.001117447496703505
$ws.Range("E5").Value = 0.05720210856645425
$ws.Range("F5").Value = 0.4538423500840025
$ws.Range("G5").Value = 0.3349609454740374
$ws.Range("H5").Value = 0.4112601324958689
$ws.Range("M5").Value = 0.6877648806897412
$ws.Range("N5").Value = 1.49843577548404
$ws.Range("O5").Value = 1.42303945874869
$ws.Range("C6").Value = 0.03963749364594094
$ws.Range("D6").Value = 0.001109459715097572
$ws.Range("E6").Value = 0.05735687108794019
$ws.Range("F6").Value = 0.4526383833542127
$ws.Range("G6").Value = 0.3336697189851634
$ws.Range("H6").Value = 0.4108576011172858
$ws.Range("M6").Value = 0.6832732823300773
$ws.Range("N6").Value = 1.492191195483997
$ws.Range("O6").Value = 1.419441944818033
$ws.Range("C7").Value = 0.04127938586999846
$ws.Range("D7").Value = 0.001164871035888027
$ws.Range("E7").Value = 0.05629570366196179
$ws.Range("F7").Value = 0.4610271394244592
$ws.Range("G7").Value = 0.342655913838243
$ws.Range("H7").Value = 0.4136776126512416
$ws.Range("M7").Value = 0.7144364355493025
$ws.Range("N7").Value = 1.53555089848183
$ws.Range("O7").Value = 1.444515944619042
$ws.Range("C8").Value = 0.04855855672559528
$ws.Range("D8").Value = 0.001408544397047962
$ws.Range("E8").Value = 0.05193833640862211
$ws.Range("F8").Value = 0.4988382434003285
$ws.Range("G8").Value = 0.3829005436759303
$ws.Range("H8").Value = 0.4267666056344837
$ws.Range("M8").Value = 0.8516086809471091
$ws.Range("N8").Value = 1.727239951245508
$ws.Range("O8").Value = 1.557731870001561
$ws.Range("C9").Value = 0.06295881630227029
$ws.Range("D9").Value = 0.00188301898717036
$ws.Range("E9").Value = 0.04458506273144214
$ws.Range("F9").Value = 0.5758768262664091
$ws.Range("G9").Value = 0.4639921232089819
$ws.Range("H9").Value = 0.4547675485921729
$ws.Range("M9").Value = 1.119274021068847
$ws.Range("N9").Value = 2.104132684427839
$ws.Range("O9").Value = 1.78912283365554
$ws.Range("C10").Value = 0.0736205858625425
$ws.Range("D10").Value = 0.002229310950426111
$ws.Range("E10").Value = 0.03992928257694128
$ws.Range("F10").Value = 0.6342846482094444
$ws.Range("G10").Value = 0.5249482728020496
$ws.Range("H10").Value = 0.4767805981864512
$ws.Range("M10").Value = 1.315060226785022
$ws.Range("N10").Value = 2.381433539546606
$ws.Range("O10").Value = 1.964994345450179
$ws.Range("C11").Value = 0.07848985866168334
$ws.Range("D11").Value = 0.002386297562104289
$ws.Range("E11").Value = 0.03797776431191924
$ws.Range("F11").Value = 0.6612584602212621
$ws.Range("G11").Value = 0.5529900601040652
$ws.Range("H11").Value = 0.4871118650747519
$ws.Range("M11").Value = 1.403929410150184
$ws.Range("N11").Value = 2.507630204147233
$ws.Range("O11").Value = 2.046311345068602
$ws.Range("C12").Value = 0.08033655473910528
$ws.Range("D12").Value = 0.002445661343060124
$ws.Range("E12").Value = 0.03726307507434856
$ws.Range("F12").Value = 0.671531487568231
$ws.Range("G12").Value = 0.5636545334886591
$ws.Range("H12").Value = 0.4910699680696098
$ws.Range("M12").Value = 1.437552522977612
$ws.Range("N12").Value = 2.555420990854714
$ws.Range("O12").Value = 2.077295028259925
$ws.Range("C13").Value = 0.07993871000778086
$ws.Range("D13").Value = 0.002432880098087509
$ws.Range("E13").Value = 0.03741590902867187
$ws.Range("F13").Value = 0.6693163907281701
$ws.Range("G13").Value = 0.5613557063034591
$ws.Range("H13").Value = 0.4902154756935317
$ws.Range("M13").Value = 1.430312532385429
$ws.Range("N13").Value = 2.5451283341975
$ws.Range("O13").Value = 2.070613624151065
$ws.Range("C14").Value = 0.07864173089369331
$ws.Range("D14").Value = 0.00239118316311604
$ws.Range("E14").Value = 0.03791847705088625
$ws.Range("F14").Value = 0.662102450091254
$ws.Range("G14").Value = 0.5538665136045893
$ws.Range("H14").Value = 0.4874365798253564
$ws.Range("M14").Value = 1.406696209058609
$ws.Range("N14").Value = 2.511561949414272
$ws.Range("O14").Value = 2.048856557864497
$ws.Range("C15").Value = 0.07784766092331097
$ws.Range("D15").Value = 0.002365631532779133
$ws.Range("E15").Value = 0.03822949139090914
$ws.Range("F15").Value = 0.6576913552988088
$ws.Range("G15").Value = 0.5492851340097218
$ws.Range("H15").Value = 0.4857404072398026
$ws.Range("M15").Value = 1.392226609244616
$ws.Range("N15").Value = 2.4910018079201
$ws.Range("O15").Value = 2.035554621181859
$ws.Range("C16").Value = 0.0733027576706462
$ws.Range("D16").Value = 0.002219040099447511
$ws.Range("E16").Value = 0.04006019987723297
$ws.Range("F16").Value = 0.6325300198437986
$ws.Range("G16").Value = 0.5231220260897658
$ws.Range("H16").Value = 0.4761118326834151
$ws.Range("M16").Value = 1.309248349344202
$ws.Range("N16").Value = 2.373186933110503
$ws.Range("O16").Value = 1.959706664292924
$ws.Range("C17").Value = 0.0705195549099642
$ws.Range("D17").Value = 0.002128967898769929
$ws.Range("E17").Value = 0.04122617200752465
$ws.Range("F17").Value = 0.6171981794308294
$ws.Range("G17").Value = 0.507152361793004
$ws.Range("H17").Value = 0.4702864774332056
$ws.Range("M17").Value = 1.258292739170741
$ws.Range("N17").Value = 2.30092119317078
$ws.Range("O17").Value = 1.913513906478215
$ws.Range("C18").Value = 0.06892053265904963
$ws.Range("D18").Value = 0.002077109934376153
$ws.Range("E18").Value = 0.04191244979674913
$ws.Range("F18").Value = 0.6084176841589937
$ws.Range("G18").Value = 0.4979964071009988
$ws.Range("H18").Value = 0.4669657502940936
$ws.Range("M18").Value = 1.228966170413059
$ws.Range("N18").Value = 2.259360927329112
$ws.Range("O18").Value = 1.887068440558778
$ws.Range("C19").Value = 0.06837943923056855
$ws.Range("D19").Value = 0.002059543159592891
$ws.Range("E19").Value = 0.04214748664310375
$ws.Range("F19").Value = 0.6054512642119789
$ws.Range("G19").Value = 0.4949013813633485
$ws.Range("H19").Value = 0.4658465311307225
$ws.Range("M19").Value = 1.219033615444374
$ws.Range("N19").Value = 2.245290345450655
$ws.Range("O19").Value = 1.87813559077091
$ws.Range("C20").Value = 0.07081564469505963
$ws.Range("D20").Value = 0.002138561537723405
$ws.Range("E20").Value = 0.04110043072404479
$ws.Range("F20").Value = 0.6188263457869141
$ws.Range("G20").Value = 0.5088493145602513
$ws.Range("H20").Value = 0.4709035040451397
$ws.Range("M20").Value = 1.263718951894802
$ws.Range("N20").Value = 2.308613511958015
$ws.Range("O20").Value = 1.918418419794193
$ws.Range("C21").Value = 0.07902260848796061
$ws.Range("D21").Value = 0.002403432877585487
$ws.Range("E21").Value = 0.03777019789664671
$ws.Range("F21").Value = 0.6642197636153924
$ws.Range("G21").Value = 0.5560650262067952
$ws.Range("H21").Value = 0.4882515619740389
$ws.Range("M21").Value = 1.413633715564913
$ws.Range("N21").Value = 2.52142116855407
$ws.Range("O21").Value = 2.055241944010334
$ws.Range("C22").Value = 0.08440272494220835
$ws.Range("D22").Value = 0.002576051577442229
$ws.Range("E22").Value = 0.03573556309468051
$ws.Range("F22").Value = 0.6942289906527463
$ws.Range("G22").Value = 0.5871896812987529
$ws.Range("H22").Value = 0.4998570324233356
$ws.Range("M22").Value = 1.51143769152533
$ws.Range("N22").Value = 2.660517831285119
$ws.Range("O22").Value = 2.145776512420753
$ws.Range("C23").Value = 0.08152973922818774
$ws.Range("D23").Value = 0.002483968368611755
$ws.Range("E23").Value = 0.03680837993740593
$ws.Range("F23").Value = 0.6781810211586645
$ws.Range("G23").Value = 0.5705532492066538
$ws.Range("H23").Value = 0.4936384200856025
$ws.Range("M23").Value = 1.459254340543723
$ws.Range("N23").Value = 2.586279464566587
$ws.Range("O23").Value = 2.097354050181195
$ws.Range("C24").Value = 0.0706817791512151
$ws.Range("D24").Value = 0.002134224484905189
$ws.Range("E24").Value = 0.04115722871708893
$ws.Range("F24").Value = 0.618090145961645
$ws.Range("G24").Value = 0.5080820437194404
$ws.Range("H24").Value = 0.4706244580260659
$ws.Range("M24").Value = 1.261265859152431
$ws.Range("N24").Value = 2.305135857049493
$ws.Range("O24").Value = 1.916200742575484
$ws.Range("C25").Value = 0.05904915692036639
$ws.Range("D25").Value = 0.001755047841117019
$ws.Range("E25").Value = 0.04644467003064112
$ws.Range("F25").Value = 0.5547216813349962
$ws.Range("G25").Value = 0.4418165940236634
$ws.Range("H25").Value = 0.4469408611618633
$ws.Range("M25").Value = 1.047011733965022
$ws.Range("N25").Value = 2.002088683901547
$ws.Range("O25").Value = 1.725505971028326
